# Misc card changes: rework "Challenge"/alt-upgrade cards, tweak recovery
# amounts on "再生" and replace "驼兽" with the new "强韧" card.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - "等级3" card: "重抽本牌" -> "从升级牌堆中选1张《等级2》替换本牌"
$ws.Range("E3").Value = "进入备选升级牌区时：如果玩家未拥有《等级2》，则从升级牌堆中选1张《等级2》替换本牌。<br>`n被动：可以使用至多包含3张牌的堆叠。"

# Row 4 - "等级4" card: same wording swap, referencing "等级3"
$ws.Range("E4").Value = "进入备选升级牌区时：如果玩家未拥有《等级3》，则从升级牌堆中选1张《等级3》替换本牌。<br>被动：可以使用至多包含4张牌的堆叠。"

# Row 5 - "等级5" card: same wording swap, referencing "等级4"
$ws.Range("E5").Value = "进入备选升级牌区时：如果玩家未拥有《等级4》，则从升级牌堆中选1张《等级4》替换本牌。<br>被动：可以使用至多包含5张牌的堆叠。"

# Row 7 - "再生" card: recovery amount 2 -> 1 (both effects)
$ws.Range("E7").Value = "主动：横置1张“体质”牌，回复1生命。<br>`n主动：将1张手牌洗回主牌堆，回复1生命。"

# Row 6 - replace "驼兽" card with the new "强韧" card
$ws.Range("A6").Value = "强韧"
$ws.Range("C6").Value = 1
$ws.Range("E6").Value = "被动：最大生命值加2。<br>`n获得本牌时：回复所有生命值。"

# Row 9 - "透支" card: maxCount 5 -> 3
$ws.Range("B9").Value = 3

# Row heights grew because of the longer replacement text (wrap-text column E)
$ws.Rows(3).RowHeight = 242.25
$ws.Rows(4).RowHeight = 242.25
$ws.Rows(5).RowHeight = 242.25
$ws.Rows(6).RowHeight = 114

# Restore the on-screen selection/scroll position recorded in the saved file
$ws.Range("E7").Select()
